$d = $word.ActiveDocument

# Paragraph indices (1-based) that need strikethrough formatting applied,
# matching the "doelgroep < OG:" ... "kale sommen ipv verhaaltjessommen"
# block, plus the "stopmotion prototyping (userscenario's)" paragraph.
$indices = @(7, 8, 9, 10, 11, 12, 13, 14, 15, 18)

foreach ($i in $indices) {
    $p = $d.Paragraphs.Item($i)
    $p.Range.Font.StrikeThrough = 1
}
